{"js": "// Apply the \"Add data structure\" edit to Data Flow/Packet.docx\n//\n// 1. Table cell \"Thang Duong Chi\" -> split into \"Thang\" (wrapped in\n//    proofErr spellStart/spellEnd) + \" Duong Chi\".\n// 2. \"Contains packets from network.\" paragraph -> wrap the run in\n//    proofErr gramStart/gramEnd.\n// 3. \"Packet\" paragraph gains the data-structure formula text\n//    (\" = Version + Protocol + Src_IP + Dest_IP + Data \") with\n//    proofErr spellStart/spellEnd around \"Src_IP\" / \"Dest_IP\", and the\n//    \"_GoBack\" bookmark moves here (from the former \"Unresolved Issue\"\n//    paragraph).\n// 4. \"1/ms\" paragraph -> split into \"1/\" + \"ms\" (wrapped in proofErr\n//    spellStart/spellEnd).\n// 5. Final paragraph: the two runs \"The type of the data flow and the\n//    volume/time \" / \"may not be correct\" (separated by the _GoBack\n//    bookmark) collapse into a single run/sentence, bookmark removed.\n\nconst PKG_OPEN =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>';\nconst PKG_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrapParagraph(innerXml) {\n  return PKG_OPEN + innerXml + PKG_CLOSE;\n}\n\nasync function replaceParagraphOoxml(paragraph, innerXml) {\n  const range = paragraph.getRange();\n  range.insertOoxml(wrapParagraph(innerXml), Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfunction findIndex(predicate, fromIndex) {\n  const start = fromIndex || 0;\n  for (let i = start; i < items.length; i++) {\n    if (predicate(items[i].text)) {\n      return i;\n    }\n  }\n  throw new Error(\"Paragraph not found\");\n}\n\n// --- 1. \"Thang Duong Chi\" table cell -----------------------------------\nconst thangIdx = findIndex((t) => t === \"Thang Duong Chi\");\nawait replaceParagraphOoxml(\n  items[thangIdx],\n  '<w:p w:rsidR=\"00C64291\" w:rsidRDefault=\"004275D2\" w:rsidP=\"00B5208A\">' +\n    '<w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Thang</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Duong Chi</w:t></w:r>' +\n    \"</w:p>\"\n);\n\n// --- 2. \"Contains packets from network.\" --------------------------------\nconst containsIdx = findIndex((t) => t === \"Contains packets from network.\");\nawait replaceParagraphOoxml(\n  items[containsIdx],\n  '<w:p w:rsidR=\"00C64291\" w:rsidRPr=\"00C64291\" w:rsidRDefault=\"004275D2\" w:rsidP=\"00C64291\">' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Contains packets from network.</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    \"</w:p>\"\n);\n\n// --- 3. \"Packet\" paragraph gains the data-structure formula -------------\nconst packetIdx = findIndex((t) => t === \"Packet\");\nawait replaceParagraphOoxml(\n  items[packetIdx],\n  '<w:p w:rsidR=\"004A3EA4\" w:rsidRPr=\"004A3EA4\" w:rsidRDefault=\"004275D2\" w:rsidP=\"004A3EA4\">' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Packet</w:t></w:r>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> = Version + Protocol + </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Src_IP</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> + </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Dest_IP</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> + Data </w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    \"</w:p>\"\n);\n\n// --- 4. \"1/ms\" -----------------------------------------------------------\nconst oneMsIdx = findIndex((t) => t === \"1/ms\");\nawait replaceParagraphOoxml(\n  items[oneMsIdx],\n  '<w:p w:rsidR=\"004A3EA4\" w:rsidRPr=\"004A3EA4\" w:rsidRDefault=\"004275D2\" w:rsidP=\"004A3EA4\">' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>1/</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>ms</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"</w:p>\"\n);\n\n// --- 5. Final \"Unresolved Issue\" paragraph: merge runs, drop bookmark ---\nconst finalIdx = findIndex(\n  (t) => t === \"The type of the data flow and the volume/time may not be correct\"\n);\nawait replaceParagraphOoxml(\n  items[finalIdx],\n  '<w:p w:rsidR=\"004A3EA4\" w:rsidRPr=\"004A3EA4\" w:rsidRDefault=\"004275D2\" w:rsidP=\"004A3EA4\">' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>The type of the data flow and the volume/time may not be correct</w:t></w:r>' +\n    \"</w:p>\"\n);\n\nawait context.sync();\n", "ps1": "# Apply the \"Add data structure\" edit to Data Flow/Packet.docx\n#\n# 1. Table cell \"Thang Duong Chi\" -> split into \"Thang\" (wrapped in\n#    proofErr spellStart/spellEnd) + \" Duong Chi\".\n# 2. \"Contains packets from network.\" paragraph -> wrap the run in\n#    proofErr gramStart/gramEnd.\n# 3. \"Packet\" paragraph gains the data-structure formula text\n#    (\" = Version + Protocol + Src_IP + Dest_IP + Data \") with\n#    proofErr spellStart/spellEnd around \"Src_IP\" / \"Dest_IP\", and the\n#    \"_GoBack\" bookmark moves here (from the former \"Unresolved Issue\"\n#    paragraph).\n# 4. \"1/ms\" paragraph -> split into \"1/\" + \"ms\" (wrapped in proofErr\n#    spellStart/spellEnd).\n# 5. Final paragraph: the two runs \"The type of the data flow and the\n#    volume/time \" / \"may not be correct\" (separated by the _GoBack\n#    bookmark) collapse into a single run/sentence, bookmark removed.\n\n$d = $word.ActiveDocument\n\n$PKG_OPEN = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$PKG_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\nfunction Get-ParaKey($text) {\n    # Paragraph text from the Word object model always carries a\n    # trailing control character (CR for normal paragraphs, CR+cell-mark\n    # for table-cell paragraphs) - strip it so comparisons are exact.\n    $s = $text\n    while ($s.Length -gt 0 -and [int][char]$s[$s.Length - 1] -lt 32) {\n        $s = $s.Substring(0, $s.Length - 1)\n    }\n    return $s\n}\n\nfunction Find-ParagraphByText($doc, $targetText) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ((Get-ParaKey $p.Range.Text) -eq $targetText) {\n            return $p\n        }\n    }\n    throw \"Paragraph not found: $targetText\"\n}\n\nfunction Set-ParagraphXml($doc, $targetText, $innerXml) {\n    $p = Find-ParagraphByText $doc $targetText\n    $r = $p.Range\n    # Quirk workaround: replacing the XML of the range that reaches all\n    # the way to the end of the document body (i.e. the very last\n    # paragraph, right before sectPr) leaves a stray empty paragraph\n    # behind. Excluding the trailing paragraph mark from the range\n    # avoids it, while still replacing all of the paragraph's content.\n    if ($r.End -eq $doc.Content.End) {\n        $r = $doc.Range($r.Start, $r.End - 1)\n    }\n    $r.InsertXML($PKG_OPEN + $innerXml + $PKG_CLOSE) | Out-Null\n}\n\n# --- 1. \"Thang Duong Chi\" table cell -------------------------------------\nSet-ParagraphXml $d \"Thang Duong Chi\" (\n    '<w:p w:rsidR=\"00C64291\" w:rsidRDefault=\"004275D2\" w:rsidP=\"00B5208A\">' +\n    '<w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Thang</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Duong Chi</w:t></w:r>' +\n    '</w:p>'\n)\n\n# --- 2. \"Contains packets from network.\" ----------------------------------\nSet-ParagraphXml $d \"Contains packets from network.\" (\n    '<w:p w:rsidR=\"00C64291\" w:rsidRPr=\"00C64291\" w:rsidRDefault=\"004275D2\" w:rsidP=\"00C64291\">' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Contains packets from network.</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '</w:p>'\n)\n\n# --- 3. \"Packet\" paragraph gains the data-structure formula ---------------\nSet-ParagraphXml $d \"Packet\" (\n    '<w:p w:rsidR=\"004A3EA4\" w:rsidRPr=\"004A3EA4\" w:rsidRDefault=\"004275D2\" w:rsidP=\"004A3EA4\">' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Packet</w:t></w:r>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> = Version + Protocol + </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Src_IP</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> + </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Dest_IP</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> + Data </w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\n)\n\n# --- 4. \"1/ms\" --------------------------------------------------------------\nSet-ParagraphXml $d \"1/ms\" (\n    '<w:p w:rsidR=\"004A3EA4\" w:rsidRPr=\"004A3EA4\" w:rsidRDefault=\"004275D2\" w:rsidP=\"004A3EA4\">' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>1/</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>ms</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>'\n)\n\n# --- 5. Final \"Unresolved Issue\" paragraph: merge runs, drop bookmark -----\nSet-ParagraphXml $d \"The type of the data flow and the volume/time may not be correct\" (\n    '<w:p w:rsidR=\"004A3EA4\" w:rsidRPr=\"004A3EA4\" w:rsidRDefault=\"004275D2\" w:rsidP=\"004A3EA4\">' +\n    '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>The type of the data flow and the volume/time may not be correct</w:t></w:r>' +\n    '</w:p>'\n)\n"}
